$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update performance metrics with new equation values
$ws.Range("B2").Value = 7500.70588235294
$ws.Range("C2").Value = 16278.9

$ws.Range("B3").Value = 161118284.617647
$ws.Range("C3").Value = 288618406.5

$ws.Range("B4").Value = 5958.02919462916
$ws.Range("C4").Value = 15444.343346087

$ws.Range("B5").Value = 88.2352941176471
$ws.Range("C5").Value = 20
